$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, pushing existing rows 11-16 down to 12-17
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with the new price record
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(11, 3).Value = "La Araucanía"
$ws.Cells.Item(11, 4).Value = 45086
$ws.Cells.Item(11, 5).Value = 9
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100104
$ws.Cells.Item(11, 8).Value = "Frutos de pepita"
$ws.Cells.Item(11, 9).Value = 100104005
$ws.Cells.Item(11, 10).Value = "Pera asiática"
$ws.Cells.Item(11, 11).Value = "Hosui"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 80
$ws.Cells.Item(11, 14).Value = 16000
$ws.Cells.Item(11, 15).Value = 16000
$ws.Cells.Item(11, 16).Value = 16000
$ws.Cells.Item(11, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(11, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(11, 19).Value = 889
$ws.Cells.Item(11, 20).Value = 18
